$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1986.4459
$ws.Range("I15").Value = 1986.4459
$ws.Range("K15").Value = 5959.3377
$ws.Range("M15").Value = -5790.3377

$ws.Range("H17").Value = 6769.0713
$ws.Range("J17").Value = 7143.6924
$ws.Range("L17").Value = 21431.0772
$ws.Range("N17").Value = -21767.0772

$ws.Range("H33").Value = 267.66666
$ws.Range("I33").Value = 267.66666
$ws.Range("K33").Value = 267.66666
$ws.Range("M33").Value = -38.66665999999998

$ws.Range("H74").Value = 17861400
$ws.Range("I74").Value = 3200
$ws.Range("J74").Value = 20837766
$ws.Range("K74").Value = 3200
$ws.Range("L74").Value = 20837766
$ws.Range("M74").Value = -2264
$ws.Range("N74").Value = -20839638

$ws.Range("H77").Value = 17861400
$ws.Range("I77").Value = 3200
$ws.Range("J77").Value = 20837766
$ws.Range("K77").Value = 16000
$ws.Range("L77").Value = 104188830
$ws.Range("M77").Value = -11320
$ws.Range("N77").Value = -104198190

$ws.Range("H112").Value = 2711192.8
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2711192.8
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 8133578.399999999
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -8135794.399999999

$ws.Range("H138").Value = 31253306
$ws.Range("I138").Value = 71430380
$ws.Range("J138").Value = 4467.778
$ws.Range("K138").Value = 214291140
$ws.Range("L138").Value = 13403.334
$ws.Range("M138").Value = -214286000
$ws.Range("N138").Value = -23683.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4372.1294
$ws.Range("I32").Value = 3919.0852
$ws.Range("J32").Value = 7414
$ws.Range("K32").Value = 3919.0852
$ws.Range("L32").Value = 7414
$ws.Range("M32").Value = -3632.0852
$ws.Range("N32").Value = -7988

$ws.Range("H61").Value = 3572.4167
$ws.Range("I61").Value = 3423.3914
$ws.Range("K61").Value = 3423.3914
$ws.Range("M61").Value = -3211.3914

$ws.Range("H110").Value = 850
$ws.Range("I110").Value = 818.8889
$ws.Range("K110").Value = 818.8889
$ws.Range("M110").Value = 1226.1111

$ws.Range("H122").Value = 2950.52
$ws.Range("I122").Value = 2338.25
$ws.Range("J122").Value = 5399.6
$ws.Range("K122").Value = 7014.75
$ws.Range("L122").Value = 16198.8
$ws.Range("M122").Value = -4564.75
$ws.Range("N122").Value = -21098.8

$ws.Range("H132").Value = 24735.262
$ws.Range("I132").Value = 2830.5334
$ws.Range("J132").Value = 65806.625
$ws.Range("K132").Value = 8491.600199999999
$ws.Range("L132").Value = 197419.875
$ws.Range("M132").Value = -5961.600199999999
$ws.Range("N132").Value = -202479.875

$ws.Range("H136").Value = 3572.4167
$ws.Range("I136").Value = 3423.3914
$ws.Range("K136").Value = 10270.1742
$ws.Range("M136").Value = -7720.174199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5088.3687
$ws.Range("I134").Value = 5259.9443
$ws.Range("K134").Value = 15779.8329
$ws.Range("M134").Value = -13244.8329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 568.5
$ws.Range("J22").Value = 650.6667
$ws.Range("L22").Value = 650.6667
$ws.Range("N22").Value = -1350.6667

$ws.Range("H52").Value = 39950
$ws.Range("J52").Value = 39950
$ws.Range("L52").Value = 39950
$ws.Range("N52").Value = -40538

$ws.Range("H58").Value = 25749.238
$ws.Range("I58").Value = 2113.25
$ws.Range("J58").Value = 40294.46
$ws.Range("K58").Value = 2113.25
$ws.Range("L58").Value = 40294.46
$ws.Range("M58").Value = -1910.25
$ws.Range("N58").Value = -40700.46

$ws.Range("H88").Value = 19671.5
$ws.Range("J88").Value = 19671.5
$ws.Range("L88").Value = 19671.5
$ws.Range("N88").Value = -20483.5

$ws.Range("H91").Value = 19671.5
$ws.Range("J91").Value = 19671.5
$ws.Range("L91").Value = 19671.5
$ws.Range("N91").Value = -22479.5

$ws.Range("H105").Value = 850.1111
$ws.Range("I105").Value = 392
$ws.Range("J105").Value = 1422.75
$ws.Range("K105").Value = 392
$ws.Range("L105").Value = 1422.75
$ws.Range("M105").Value = 1355
$ws.Range("N105").Value = -4916.75

$ws.Range("H132").Value = 3468
$ws.Range("I132").Value = 1843.8334
$ws.Range("J132").Value = 8340.5
$ws.Range("K132").Value = 5531.5002
$ws.Range("L132").Value = 25021.5
$ws.Range("M132").Value = -3001.5002
$ws.Range("N132").Value = -30081.5

$ws.Range("H134").Value = 1288.6111
$ws.Range("I134").Value = 1141.5
$ws.Range("J134").Value = 1803.5
$ws.Range("K134").Value = 3424.5
$ws.Range("L134").Value = 5410.5
$ws.Range("M134").Value = -889.5
$ws.Range("N134").Value = -10480.5

$ws.Range("H136").Value = 25749.238
$ws.Range("I136").Value = 2113.25
$ws.Range("J136").Value = 40294.46
$ws.Range("K136").Value = 6339.75
$ws.Range("L136").Value = 120883.38
$ws.Range("M136").Value = -3789.75
$ws.Range("N136").Value = -125983.38

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 699.84
$ws.Range("J131").Value = 717.35486
$ws.Range("L131").Value = 2152.06458
$ws.Range("N131").Value = -12232.06458

$ws.Range("H140").Value = 2735
$ws.Range("I140").Value = 1369.091
$ws.Range("J140").Value = 3890.7693
$ws.Range("K140").Value = 4107.272999999999
$ws.Range("L140").Value = 11672.3079
$ws.Range("M140").Value = 1072.727000000001
$ws.Range("N140").Value = -22032.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 4117764.8
$ws.Range("I7").Value = 5000000
$ws.Range("J7").Value = 3125250
$ws.Range("K7").Value = 5000000
$ws.Range("L7").Value = 3125250
$ws.Range("M7").Value = -4999888
$ws.Range("N7").Value = -3125474

$ws.Range("H8").Value = 4117764.8
$ws.Range("I8").Value = 5000000
$ws.Range("J8").Value = 3125250
$ws.Range("K8").Value = 5000000
$ws.Range("L8").Value = 3125250
$ws.Range("M8").Value = -4999861
$ws.Range("N8").Value = -3125528

$ws.Range("H11").Value = 7901600.5
$ws.Range("I11").Value = 10142857
$ws.Range("J11").Value = 2672001.2
$ws.Range("K11").Value = 10142857
$ws.Range("L11").Value = 2672001.2
$ws.Range("M11").Value = -10142718
$ws.Range("N11").Value = -2672279.2

$ws.Range("H12").Value = 6593333.5
$ws.Range("I12").Value = 6593333.5
$ws.Range("K12").Value = 6593333.5
$ws.Range("M12").Value = -6593193.5

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 725.5454999999999
$ws.Range("I16").Value = 678.1
$ws.Range("K16").Value = 678.1
$ws.Range("M16").Value = -508.1

$ws.Range("H22").Value = 5087.5
$ws.Range("I22").Value = 3450.3333
$ws.Range("J22").Value = 9999
$ws.Range("K22").Value = 3450.3333
$ws.Range("L22").Value = 9999
$ws.Range("M22").Value = -3155.3333
$ws.Range("N22").Value = -10589

$ws.Range("H27").Value = 5087.5
$ws.Range("I27").Value = 3450.3333
$ws.Range("J27").Value = 9999
$ws.Range("K27").Value = 3450.3333
$ws.Range("L27").Value = 9999
$ws.Range("M27").Value = -3343.3333
$ws.Range("N27").Value = -10213

$ws.Range("H136").Value = 1794.3158
$ws.Range("I136").Value = 1613.7142
$ws.Range("K136").Value = 4841.142599999999
$ws.Range("M136").Value = -2291.142599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 22442786
$ws.Range("I136").Value = 28674884
$ws.Range("K136").Value = 86024652
$ws.Range("M136").Value = -86022102
